$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort the data rows (A2:D8) in ascending order of column A ("time (s)"),
# matching the re-ordering performed during needle calibration.
$range = $ws.Range("A2:D8")
$key1 = $ws.Range("A2:A8")

# Sort(Key1, Order1, Key2, Type, Order2, Key3, Order3, Header, OrderCustom, MatchCase, Orientation, SortMethod)
# Order1 = 1 (xlAscending), Header = 2 (xlNo - range has no header row)
$range.Sort($key1, 1, $null, $null, 1, $null, 1, 2, $false, $null, 1, 1)
